$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 86; $row++) {
    $ws.Cells.Item($row, 15).Value = "2022-07-29 20:58:22"
}
